$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'278.58"
$ws.Range("E2").Value = "'1.11%"
$ws.Range("D3").Value = "'27.40"
$ws.Range("E3").Value = "'0.08%"
$ws.Range("D4").Value = "'4.817"
$ws.Range("E4").Value = "'0.24%"
$ws.Range("D5").Value = "'0.06378"
$ws.Range("E5").Value = "'0.43%"
$ws.Range("D6").Value = "'7.037"
$ws.Range("E6").Value = "'1.19%"
$ws.Range("D7").Value = "'1.295"
$ws.Range("E7").Value = "'-3.75%"
$ws.Range("D8").Value = "'0.8920"
$ws.Range("E8").Value = "'1.48%"
$ws.Range("D9").Value = "'0.1519"
$ws.Range("E9").Value = "'-0.14%"
$ws.Range("D10").Value = "'0.05656"
$ws.Range("E10").Value = "'11.67%"
$ws.Range("D11").Value = "'0.07454"
$ws.Range("E11").Value = "'-0.77%"
$ws.Range("D12").Value = "'0.02909"
$ws.Range("E12").Value = "'-2.06%"
$ws.Range("D13").Value = "'0.08973"
$ws.Range("E13").Value = "'-0.64%"
$ws.Range("D14").Value = "'0.001573"
$ws.Range("E14").Value = "'0.30%"
$ws.Range("D15").Value = "'0.0006365"
$ws.Range("E15").Value = "'-0.55%"
$ws.Range("D16").Value = "'0.006128"
$ws.Range("E16").Value = "'7.69%"
$ws.Range("D17").Value = "'3.471"
$ws.Range("D18").Value = "'3.300"
$ws.Range("E18").Value = "'-0.02%"
$ws.Range("D19").Value = "'2.232"
$ws.Range("E19").Value = "'-2.29%"
$ws.Range("E21").Value = "'-0.30%"
$ws.Range("D22").Value = "'3.908"
$ws.Range("E22").Value = "'-0.29%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.1504"
$ws.Range("E23").Value = "'9.01%"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").Value = "'0.04382"
$ws.Range("E24").Value = "'-0.71%"
$ws.Range("D25").Value = "'0.001174"
$ws.Range("E25").Value = "'0.29%"
$ws.Range("D26").Value = "'0.004280"
$ws.Range("E26").Value = "'10.78%"
$ws.Range("D28").Value = "'0.0001178"
$ws.Range("E28").Value = "'-1.80%"
$ws.Range("D29").Value = "'0.0001651"
$ws.Range("E29").Value = "'-14.74%"
$ws.Range("D40").Value = "'0.04029"
$ws.Range("E40").Value = "'-3.82%"
$ws.Range("D41").Value = "'0.006737"
$ws.Range("E41").Value = "'-1.71%"
$ws.Range("D42").Value = "'0.1393"
$ws.Range("E42").Value = "'18.07%"
$ws.Range("D43").Value = "'0.002037"
$ws.Range("E43").Value = "'0.85%"
$ws.Range("D44").Value = "'0.01116"
$ws.Range("E44").Value = "'-3.34%"
$ws.Range("D45").Value = "'0.00005524"
$ws.Range("E45").Value = "'6.98%"
$ws.Range("D46").Value = "'1.628"
$ws.Range("E46").Value = "'9.31%"
$ws.Range("D47").Value = "'0.01847"
$ws.Range("E47").Value = "'-19.71%"
